# before.xlsx -> after.xlsx
# Commit message: "3 white 3 flash 4 gray done by Jim"
#
# Appends a new block of 20 rows (PartipantID = 3) to the BarChart data
# table on Sheet1, mirroring the existing PartipantID = 5 / 4 blocks that
# are already on the sheet, and updates the sheet active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: PartipantID, Index, Vis, VisID, Error, TruePerc, ReportPerc
$newRows = @(
    @(3, 0, "BarChart", 0, 2.8803084000000001, 32.238075000000002, 25),
    @(3, 1, "BarChart", 0, 2.2935064000000001, 59.777462, 55),
    @(3, 2, "BarChart", 0, 1.7930874000000001, 33.340556999999997, 30),
    @(3, 3, "BarChart", 0, 1.6327113, 72.024050000000003, 75),
    @(3, 4, "BarChart", 0, 1.220572, 37.794609999999999, 40),
    @(3, 5, "BarChart", 0, 2.3356661999999999, 40.077159999999999, 45),
    @(3, 6, "BarChart", 0, 2.6220705999999998, 43.968670000000003, 50),
    @(3, 7, "BarChart", 0, -1.7377951, 33.174827999999998, 33),
    @(3, 8, "BarChart", 0, 2.0890390000000001, 75.870350000000002, 80),
    @(3, 9, "BarChart", 0, 0.66562736, 76.461259999999996, 75),
    @(3, 10, "BarChart", 0, 1.1047522999999999, 62.025620000000004, 60),
    @(3, 11, "BarChart", 0, -0.85019330000000004, 80.42971, 80),
    @(3, 12, "BarChart", 0, 3.0282588000000001, 28.033245000000001, 20),
    @(3, 13, "BarChart", 0, 1.5472918, 32.79768, 30),
    @(3, 14, "BarChart", 0, -0.37640423000000001, 24.354645000000001, 25),
    @(3, 15, "BarChart", 0, 1.9364768000000001, 43.702697999999998, 40),
    @(3, 16, "BarChart", 0, 0.67729985999999998, 96.474143999999995, 95),
    @(3, 17, "BarChart", 0, 1.3410127000000001, 62.408290000000001, 60),
    @(3, 18, "BarChart", 0, 1.240032, 87.237039999999993, 85),
    @(3, 19, "BarChart", 0, 0.23042393, 73.951819999999998, 75)
)

$startRow = 42
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    for ($col = 1; $col -le 7; $col++) {
        $ws.Cells.Item($r, $col).Value = $rowData[$col - 1]
    }
}

# Update the active selection to match the saved sheet view state
$ws.Range("E15").Select()
